$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1428.3684
$ws.Range("I33").Value = 1042.75
$ws.Range("J33").Value = 2089.4285
$ws.Range("K33").Value = 1042.75
$ws.Range("L33").Value = 2089.4285
$ws.Range("M33").Value = -813.75
$ws.Range("N33").Value = -2547.4285

$ws.Range("H40").Value = 3928.5
$ws.Range("J40").Value = 9499.75
$ws.Range("L40").Value = 9499.75
$ws.Range("N40").Value = -9849.75

$ws.Range("H41").Value = 820.619
$ws.Range("I41").Value = 845.8
$ws.Range("K41").Value = 845.8
$ws.Range("M41").Value = -405.8

$ws.Range("H55").Value = 84.47619
$ws.Range("I55").Value = 84.3125
$ws.Range("K55").Value = 84.3125
$ws.Range("M55").Value = 129.6875

$ws.Range("H103").Value = 1133.1765
$ws.Range("I103").Value = 232.66667
$ws.Range("J103").Value = 1624.3636
$ws.Range("K103").Value = 698.00001
$ws.Range("L103").Value = 4873.0908
$ws.Range("M103").Value = -112.00001
$ws.Range("N103").Value = -6045.0908

$ws.Range("H106").Value = 6949.4
$ws.Range("I106").Value = 8321.333000000001
$ws.Range("J106").Value = 4891.5
$ws.Range("K106").Value = 8321.333000000001
$ws.Range("L106").Value = 4891.5
$ws.Range("M106").Value = -7690.333000000001
$ws.Range("N106").Value = -6153.5

$ws.Range("H112").Value = 1732.4814
$ws.Range("J112").Value = 1949.35
$ws.Range("L112").Value = 5848.049999999999
$ws.Range("N112").Value = -8064.049999999999

$ws.Range("H137").Value = 1242.6666
$ws.Range("I137").Value = 1236.1666
$ws.Range("K137").Value = 3708.4998
$ws.Range("M137").Value = -1158.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2349.5312
$ws.Range("I74").Value = 1321.6818
$ws.Range("K74").Value = 1321.6818
$ws.Range("M74").Value = -447.6818000000001

$ws.Range("H77").Value = 2349.5312
$ws.Range("I77").Value = 1321.6818
$ws.Range("K77").Value = 6608.409000000001
$ws.Range("M77").Value = -2240.409000000001

$ws.Range("H132").Value = 3177.6667
$ws.Range("I132").Value = 2570.7646
$ws.Range("K132").Value = 7712.293799999999
$ws.Range("M132").Value = -5182.293799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1169.1052
$ws.Range("I86").Value = 1217.125
$ws.Range("J86").Value = 1134.1818
$ws.Range("K86").Value = 1217.125
$ws.Range("L86").Value = 1134.1818
$ws.Range("M86").Value = -94.125
$ws.Range("N86").Value = -3380.1818

$ws.Range("H89").Value = 1169.1052
$ws.Range("I89").Value = 1217.125
$ws.Range("J89").Value = 1134.1818
$ws.Range("K89").Value = 6085.625
$ws.Range("L89").Value = 5670.909000000001
$ws.Range("M89").Value = -469.625
$ws.Range("N89").Value = -16902.909

$ws.Range("H94").Value = 6251105
$ws.Range("I94").Value = 8929497
$ws.Range("K94").Value = 8929497
$ws.Range("M94").Value = -8929046

$ws.Range("H107").Value = 66670
$ws.Range("J107").Value = 3608.2
$ws.Range("L107").Value = 3608.2
$ws.Range("N107").Value = -7448.2

$ws.Range("H126").Value = 69650
$ws.Range("J126").Value = 69650
$ws.Range("L126").Value = 69650
$ws.Range("N126").Value = -79530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1608.2593
$ws.Range("I31").Value = 1087.2667
$ws.Range("J31").Value = 2259.5
$ws.Range("K31").Value = 1087.2667
$ws.Range("L31").Value = 2259.5
$ws.Range("M31").Value = -792.2666999999999
$ws.Range("N31").Value = -2849.5

$ws.Range("H34").Value = 1608.2593
$ws.Range("I34").Value = 1087.2667
$ws.Range("J34").Value = 2259.5
$ws.Range("K34").Value = 1087.2667
$ws.Range("L34").Value = 2259.5
$ws.Range("M34").Value = -885.2666999999999
$ws.Range("N34").Value = -2663.5

$ws.Range("H75").Value = 58428.57
$ws.Range("J75").Value = 58428.57
$ws.Range("L75").Value = 58428.57
$ws.Range("N75").Value = -60424.57

$ws.Range("H78").Value = 58428.57
$ws.Range("J78").Value = 58428.57
$ws.Range("L78").Value = 175285.71
$ws.Range("N78").Value = -185269.71

$ws.Range("H107").Value = 1373.9333
$ws.Range("J107").Value = 2271
$ws.Range("L107").Value = 2271
$ws.Range("N107").Value = -6111

$ws.Range("H134").Value = 37038370
$ws.Range("I134").Value = 37038370
$ws.Range("K134").Value = 111115110
$ws.Range("M134").Value = -111112575

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 958.3333
$ws.Range("J68").Value = 1112.5
$ws.Range("L68").Value = 3337.5
$ws.Range("N68").Value = -4959.5

$ws.Range("H71").Value = 958.3333
$ws.Range("J71").Value = 1112.5
$ws.Range("L71").Value = 10012.5
$ws.Range("N71").Value = -18124.5

$ws.Range("H76").Value = 6814.6665
$ws.Range("J76").Value = 6814.6665
$ws.Range("L76").Value = 20443.9995
$ws.Range("N76").Value = -21209.9995

$ws.Range("H79").Value = 6814.6665
$ws.Range("J79").Value = 6814.6665
$ws.Range("L79").Value = 20443.9995
$ws.Range("N79").Value = -23095.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H113").Value = 28524.25
$ws.Range("I113").Value = 28524.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 28524.25
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -26354.25
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 21800

$ws.Range("H20").Value = 1000
$ws.Range("I20").Value = 600
$ws.Range("J20").Value = 1500
$ws.Range("K20").Value = 600
$ws.Range("L20").Value = 1500
$ws.Range("M20").Value = -374
$ws.Range("N20").Value = -1952

$ws.Range("H22").Value = 6583.3335
$ws.Range("J22").Value = 10600
$ws.Range("L22").Value = 10600
$ws.Range("N22").Value = -11190

$ws.Range("H27").Value = 6583.3335
$ws.Range("J27").Value = 10600
$ws.Range("L27").Value = 10600
$ws.Range("N27").Value = -10814

$ws.Range("H40").Value = 4748.5
$ws.Range("I40").Value = 4914.1113
$ws.Range("J40").Value = 4535.5713
$ws.Range("K40").Value = 4914.1113
$ws.Range("L40").Value = 4535.5713
$ws.Range("M40").Value = -4778.1113
$ws.Range("N40").Value = -4807.5713

$ws.Range("H55").Value = 1748.7778
$ws.Range("I55").Value = 481.06668
$ws.Range("J55").Value = 3333.4167
$ws.Range("K55").Value = 481.06668
$ws.Range("L55").Value = 3333.4167
$ws.Range("M55").Value = -308.06668
$ws.Range("N55").Value = -3679.4167

$ws.Range("H136").Value = 3687.111
$ws.Range("I136").Value = 3338.8
$ws.Range("K136").Value = 10016.4
$ws.Range("M136").Value = -7466.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 13003.728
$ws.Range("I15").Value = 4000
$ws.Range("J15").Value = 13904.1
$ws.Range("K15").Value = 4000
$ws.Range("L15").Value = 13904.1
$ws.Range("N15").Value = -14480.1
$ws.Range("M15").Value = -3712

$ws.Range("H104").Value = 45099.6
$ws.Range("J104").Value = 45099.6
$ws.Range("L104").Value = 45099.6
$ws.Range("N104").Value = -52087.6

$ws.Range("H107").Value = 15905.25
$ws.Range("I107").Value = 3648.5
$ws.Range("J107").Value = 36333.168
$ws.Range("K107").Value = 10945.5
$ws.Range("L107").Value = 108999.504
$ws.Range("M107").Value = -9025.5
$ws.Range("N107").Value = -112839.504
